# Applies the "test_requests.xlsx" update:
#  - adds 12 new header columns (AC1:AN1) with the same bold/centered/
#    bordered header style used by the rest of row 1
#  - converts M6 / N6 / AA6 from text to real numbers
#  - appends a brand new row 7 ("Images" instance) with values spread
#    across the old and new columns (the remaining row 7 cells stay
#    blank, matching the blank cells already used throughout the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header row cells AC1:AN1
# ---------------------------------------------------------------------

# Carry over the bold / centered / bordered header formatting used by the
# rest of row 1 onto the new header cells before filling in their text.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AN1").PasteSpecial(-4122)

$headers = @{
    "AC1" = "Polarisers Type"
    "AD1" = "Polarisers Lamination"
    "AE1" = "Number of Polarisers"
    "AF1" = "Orientation of Pol1"
    "AG1" = "Orientation of Cell Alignment Axis"
    "AH1" = "Orientation of Pol2"
    "AI1" = "Voltage Range"
    "AJ1" = "Voltage Single Point"
    "AK1" = "Voltage Sweep"
    "AL1" = "Tool Setup"
    "AM1" = "Tool Angle of Incidence"
    "AN1" = "Sample Number"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# ---------------------------------------------------------------------
# 2. Row 6: M6, N6, AA6 become real numbers instead of text
# ---------------------------------------------------------------------
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 0
$ws.Range("AA6").Value = 50

# ---------------------------------------------------------------------
# 3. Brand new row 7 ("Images" instance)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "Images"
$ws.Range("H7").Value = "31cdb862-d715-45fd-b3d5-d163917dec02"
$ws.Range("Z7").Value = "f4d2fe47-fdc7-4bdd-b2d8-9b526637c3ac"
$ws.Range("AL7").Value = "Nikon camera"
$ws.Range("AM7").Value = "on axis"

# These look numeric but must stay stored as text, like the other
# "(Images)" / instance columns in this sheet - write them as formulas
# that yield text, then flatten to plain values so no numeric conversion
# happens and no new number-format style gets created.
$textCells = @{
    "AC7" = "normal"
    "AD7" = "loose"
    "AE7" = "2"
    "AF7" = "0"
    "AG7" = "45"
    "AH7" = "90"
    "AI7" = "nan"
    "AJ7" = "50"
    "AK7" = "nan"
    "AN7" = "4"
}
foreach ($addr in $textCells.Keys) {
    $v = $textCells[$addr]
    $ws.Range($addr).Formula = "=""" + $v + """"
}
$ws.Range("AC7:AN7").Copy()
$ws.Range("AC7:AN7").PasteSpecial(-4163)

$ws.Range("A1").Select()
